# Score_iterations.xlsx -- "get pair with its sum closest to zero" session update
# Continues the GFG practice log: marks two existing problems "done", fixes a
# typo'd timestamp, and appends the next three problems worked on (incl. the
# start of the "Linked Lists" topic).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Mark "Majority Element" (row 120) and "Largest Sum Contiguous Subarray"
#     (row 121) as done, and fix the mis-typed semicolon in row 121's time. ---
$ws.Range("D120").Value = "done"
$ws.Range("D121").Value = "done"
$ws.Range("E121").Value = "1:22 -1:37"

# --- "Two elements whose sum is closest to zero" (row 127): mark its status
#     and correct the open-ended end time now that it's finished. ---
$ws.Range("D127").Value = "done, 2 errors"
$ws.Range("E127").Value = "5:00 - 5:18"

# --- New freeform notes row (row 128). ---
$ws.Range("D128").Value = "coding, doesn't work"
$ws.Range("E128").Value = "x - 5:22"
$ws.Range("F128").Value = "tmp"
$ws.Range("F128").ClearContents()

# --- New problem row 129: Check for Majority Element in a sorted array ---
$ws.Range("B129").Value = "GFG"
$ws.Range("C129").Value = "Check for Majority Element in a sorted array"
$ws.Range("D129").Value = "O(logn) not done, O(n) approach done"
$ws.Range("E129").Value = "5:32- 5:55"
$ws.Hyperlinks.Add($ws.Range("A129"), "https://www.geeksforgeeks.org/check-for-majority-element-in-a-sorted-array/")
$ws.Range("A129").Style = "Hyperlink"

# --- New topic row 130: start of "Linked Lists" topic, first problem ---
$ws.Range("A130").Value = "Linked Lists"
$ws.Range("A130").Style = $ws.Range("A119").Style
$ws.Range("C130").Value = "Given only a pointer/reference to a node to be deleted in a singly linked list, how do you delete it?"
$ws.Range("E130").Value = "6:04 - 6:19"
$ws.Hyperlinks.Add($ws.Range("B130"), "https://www.geeksforgeeks.org/given-only-a-pointer-to-a-node-to-be-deleted-in-a-singly-linked-list-how-do-you-delete-it/")
$ws.Range("B130").Style = "Hyperlink"

# --- New problem row 131: Detect loop in a linked list ---
$ws.Range("B131").Value = "GFG"
$ws.Range("C131").Value = "Detect loop in a linked list"
$ws.Range("E131").Value = "6:21 - x"
$ws.Hyperlinks.Add($ws.Range("A131"), "https://www.geeksforgeeks.org/detect-loop-in-a-linked-list/")
$ws.Range("A131").Style = "Hyperlink"

# --- Reflect the scrolled viewport / active selection after the edits. ---
$ws.Range("E132").Select()
